# Split the run " works on different products and customers it would be
# beneficial to be able to " into three runs by applying a green
# highlight to the middle portion: "on different products and customers
# it would be beneficial".
$d = $word.ActiveDocument

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute(
    "on different products and customers it would be beneficial",
    $true,
    $true,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "",
    0
) | Out-Null

# wdBrightGreen = 4
$rng.Font.HighlightColorIndex = 4
